$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.317380666732788
$ws.Range("B1").Value = 2.416317224502563
$ws.Range("C1").Value = 4.655238151550293
$ws.Range("D1").Value = 2.548300981521606
$ws.Range("E1").Value = 0.9577957391738892
